$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) entirely.
$ws.Rows(26).Delete()

# After the above deletion, "SC 92" (originally row 28) has shifted up to row 27.
# Delete it entirely as well.
$ws.Rows(27).Delete()

# Individual cell value corrections (imputation changes), using the new row numbers.
$ws.Range("F2").Value = $null
$ws.Range("D6").Value = -14.2
$ws.Range("D8").Value = $null
$ws.Range("D18").Value = -15.2
$ws.Range("D20").Value = $null
$ws.Range("D23").Value = -13.9
$ws.Range("D25").Value = $null

$ws.Range("C27").Value = 10
$ws.Range("C28").Value = $null
$ws.Range("C29").Value = $null
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("F30").Value = 16.89
$ws.Range("C32").Value = $null
